$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vfunction_list")

# Insert a new row above the "$set_param" row (row 9) to hold the new
# "$add_document" macro entry, shifting the existing rows down.
$ws.Rows(9).Insert()

$ws.Range("A9").Value = "`$add_document"
$ws.Range("B9").Value = "function"
$ws.Range("C9").Value = "Refer to @<bookmark:@module:add_document>"

# Make vfunction_list the active sheet/tab and select the new cell,
# matching the saved view state in the workbook.
$ws.Activate()
[void]$ws.Range("A9").Select()
